$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header E1 "EQP" -> "Flag3", add new header F1 "Flag4"
$ws.Range("E1").Value = "Flag3"
$ws.Range("F1").Value = "Flag4"

# New F1 cell should carry the same (centered) header style as the rest of row 1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data cells B2:E7 move from the (applyFont=false) centered style onto the
# (applyFont=true) centered style already used by column A / the header row,
# which collapses the redundant duplicate style.
$ws.Range("A1").Copy()
$ws.Range("B2:E7").PasteSpecial(-4122)

# Row heights: rows 1 and 3 shrink to match the other "short" rows (12.8pt)
$ws.Rows(1).RowHeight = 12.8
$ws.Rows(3).RowHeight = 12.8

# Data validation now excludes the two new flag columns (E:F) from the header row rule
$ws.Range("A1:IQ1").Validation.Delete()
$ws.Range("A2:A1007").Validation.Delete()

$v1 = $ws.Range("A1:D1").Validation
$v1.Add(3, 1, 1, '"qwer%yuiop_1234567890"', 0)
$v1.ErrorTitle = "Not Applicable"
$v1.ErrorMessage = "Cannot change the value"
$v1.IgnoreBlank = $false
$v1.ShowInput = $false
$v1.ShowError = $true

$v2 = $ws.Range("G1:IQ1").Validation
$v2.Add(3, 1, 1, '"qwer%yuiop_1234567890"', 0)
$v2.ErrorTitle = "Not Applicable"
$v2.ErrorMessage = "Cannot change the value"
$v2.IgnoreBlank = $false
$v2.ShowInput = $false
$v2.ShowError = $true

$v3 = $ws.Range("A2:A1007").Validation
$v3.Add(3, 1, 1, '"qwer%yuiop_1234567890"', 0)
$v3.ErrorTitle = "Not Applicable"
$v3.ErrorMessage = "Cannot change the value"
$v3.IgnoreBlank = $false
$v3.ShowInput = $false
$v3.ShowError = $true

# Selection moved to D5
$ws.Range("D5").Select()
